$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 - subject id headers (B:E now reflect subjects 15,16,15,16)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - updated meanEMG / legmaxROM values for subjects 15 & 16
$ws.Range("B2").Value = 88.269966626616949
$ws.Range("C2").Value = 60.780889897022746
$ws.Range("D2").Value = 48.172639143638058
$ws.Range("E2").Value = 55.357330951788185

# Row 3 - updated values; C3 no longer has data for this subject/trial
$ws.Range("B3").Value = 70.332433599997032
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 34.857935414863888
$ws.Range("E3").Value = 56.621218705819487

# Reflect the selection used while editing this range
$ws.Range("B1:E3").Select()
